$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (escape-character torture test name)
$ws.Name = "& & `" > < "

# Update cell A1 with the new escape-character torture test string
$ws.Cells.Item(1, 1).Value = "&&apos; &amp; &quot; &lt; &gt; &apos;"
